$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- RuleSet name -------------------------------------------------------
$ws.Range("D2").Value2 = "Save Complaint Rules"

# --- Rule table title ----------------------------------------------------
$ws.Range("C10").Value2 = "RuleTable Save Complaint Rules"

# --- Generalize the ACTION snippet into a templated setter call ---------
$ws.Range("D13").Value2 = "`$complaint.`$1(`$2);"

# --- Update the existing "Assign Complaint Number" rule row (15) --------
$ws.Range("D15").Value2 = "setComplaintNumber, dateFormat('yyyyMMdd') + '_' + `$complaint.getComplaintId()"

# --- Add a new rule row (16) for the Alfresco folder path ---------------
$ws.Range("B16").Value2 = "Assign Alfresco Folder"
$ws.Range("C16").Value2 = "ecmFolderId"
$ws.Range("D16").Value2 = "setEcmFolderPath, '/Sites/acm/documentLibrary/Complaints/' + dateFormat('yyyyMMdd') + '_' + `$complaint.getComplaintId()"

# --- Widen column D to fit the longer action text ------------------------
# (COM ColumnWidth is pixel-quantized; 107.2 is the closest settable value
# that round-trips to the target stored width of ~108.0357)
$ws.Columns.Item(4).ColumnWidth = 107.2

# --- Move the active selection to D17 ------------------------------------
$ws.Range("D17").Select()
